$d = $word.ActiveDocument

$d.Content.Find.Execute("67-59=8", $true, $false, $false, $false, $false, $true, 1, $false, "81-69=12", 2)
$d.Content.Find.Execute("78+21=99", $true, $false, $false, $false, $false, $true, 1, $false, "53-40=13", 2)
$d.Content.Find.Execute("74+10=84", $true, $false, $false, $false, $false, $true, 1, $false, "44+47=91", 2)
$d.Content.Find.Execute("35+1=36", $true, $false, $false, $false, $false, $true, 1, $false, "4+90=94", 2)
$d.Content.Find.Execute("76-51=25", $true, $false, $false, $false, $false, $true, 1, $false, "8+88=96", 2)
$d.Content.Find.Execute("94-8=86", $true, $false, $false, $false, $false, $true, 1, $false, "91-66=25", 2)
$d.Content.Find.Execute("42+15=57", $true, $false, $false, $false, $false, $true, 1, $false, "11+81=92", 2)
$d.Content.Find.Execute("41+27=68", $true, $false, $false, $false, $false, $true, 1, $false, "69+22=91", 2)
$d.Content.Find.Execute("70-46=24", $true, $false, $false, $false, $false, $true, 1, $false, "39-11=28", 2)
$d.Content.Find.Execute("44+48=92", $true, $false, $false, $false, $false, $true, 1, $false, "94-71=23", 2)
$d.Content.Find.Execute("51+27=78", $true, $false, $false, $false, $false, $true, 1, $false, "70+14=84", 2)
$d.Content.Find.Execute("54+12=66", $true, $false, $false, $false, $false, $true, 1, $false, "96-5=91", 2)
$d.Content.Find.Execute("87-32=55", $true, $false, $false, $false, $false, $true, 1, $false, "19-16=3", 2)
$d.Content.Find.Execute("51-34=17", $true, $false, $false, $false, $false, $true, 1, $false, "39-30=9", 2)
$d.Content.Find.Execute("2+62=64", $true, $false, $false, $false, $false, $true, 1, $false, "23+42=65", 2)
$d.Content.Find.Execute("19+27=46", $true, $false, $false, $false, $false, $true, 1, $false, "11+47=58", 2)
$d.Content.Find.Execute("9+10=19", $true, $false, $false, $false, $false, $true, 1, $false, "70+14=84", 2)
$d.Content.Find.Execute("77-47=30", $true, $false, $false, $false, $false, $true, 1, $false, "40-8=32", 2)
$d.Content.Find.Execute("55+0=55", $true, $false, $false, $false, $false, $true, 1, $false, "69+27=96", 2)
$d.Content.Find.Execute("72-54=18", $true, $false, $false, $false, $false, $true, 1, $false, "97-13=84", 2)
$d.Content.Find.Execute("4+65=69", $true, $false, $false, $false, $false, $true, 1, $false, "2+88=90", 2)
$d.Content.Find.Execute("97-31=66", $true, $false, $false, $false, $false, $true, 1, $false, "8+16=24", 2)
$d.Content.Find.Execute("41+36=77", $true, $false, $false, $false, $false, $true, 1, $false, "8+18=26", 2)
$d.Content.Find.Execute("77+18=95", $true, $false, $false, $false, $false, $true, 1, $false, "12+70=82", 2)
$d.Content.Find.Execute("30-20=10", $true, $false, $false, $false, $false, $true, 1, $false, "92-23=69", 2)
$d.Content.Find.Execute("78+9=87", $true, $false, $false, $false, $false, $true, 1, $false, "50-48=2", 2)
$d.Content.Find.Execute("18-16=2", $true, $false, $false, $false, $false, $true, 1, $false, "56-56=0", 2)
$d.Content.Find.Execute("69+26=95", $true, $false, $false, $false, $false, $true, 1, $false, "85-60=25", 2)
$d.Content.Find.Execute("77-53=24", $true, $false, $false, $false, $false, $true, 1, $false, "12+33=45", 2)
$d.Content.Find.Execute("82+12=94", $true, $false, $false, $false, $false, $true, 1, $false, "21+28=49", 2)
$d.Content.Find.Execute("95-34=61", $true, $false, $false, $false, $false, $true, 1, $false, "43+17=60", 2)
$d.Content.Find.Execute("25+74=99", $true, $false, $false, $false, $false, $true, 1, $false, "87-71=16", 2)
$d.Content.Find.Execute("58-50=8", $true, $false, $false, $false, $false, $true, 1, $false, "8+69=77", 2)
$d.Content.Find.Execute("93-22=71", $true, $false, $false, $false, $false, $true, 1, $false, "5+16=21", 2)
$d.Content.Find.Execute("11+0=11", $true, $false, $false, $false, $false, $true, 1, $false, "93-76=17", 2)
$d.Content.Find.Execute("17+61=78", $true, $false, $false, $false, $false, $true, 1, $false, "73-22=51", 2)
$d.Content.Find.Execute("84-19=65", $true, $false, $false, $false, $false, $true, 1, $false, "39-0=39", 2)
$d.Content.Find.Execute("98-30=68", $true, $false, $false, $false, $false, $true, 1, $false, "95-72=23", 2)
$d.Content.Find.Execute("71+7=78", $true, $false, $false, $false, $false, $true, 1, $false, "48-15=33", 2)
$d.Content.Find.Execute("8+65=73", $true, $false, $false, $false, $false, $true, 1, $false, "67-46=21", 2)
$d.Content.Find.Execute("82-13=69", $true, $false, $false, $false, $false, $true, 1, $false, "67-20=47", 2)
$d.Content.Find.Execute("71-48=23", $true, $false, $false, $false, $false, $true, 1, $false, "16+40=56", 2)
$d.Content.Find.Execute("28+60=88", $true, $false, $false, $false, $false, $true, 1, $false, "84-67=17", 2)
$d.Content.Find.Execute("65-52=13", $true, $false, $false, $false, $false, $true, 1, $false, "28+40=68", 2)
$d.Content.Find.Execute("70+6=76", $true, $false, $false, $false, $false, $true, 1, $false, "7+37=44", 2)
$d.Content.Find.Execute("65+18=83", $true, $false, $false, $false, $false, $true, 1, $false, "29+34=63", 2)
$d.Content.Find.Execute("27+11=38", $true, $false, $false, $false, $false, $true, 1, $false, "74-53=21", 2)
$d.Content.Find.Execute("44-0=44", $true, $false, $false, $false, $false, $true, 1, $false, "57-49=8", 2)
$d.Content.Find.Execute("11-3=8", $true, $false, $false, $false, $false, $true, 1, $false, "9+69=78", 2)
$d.Content.Find.Execute("90-85=5", $true, $false, $false, $false, $false, $true, 1, $false, "2+86=88", 2)
$d.Content.Find.Execute("4+50=54", $true, $false, $false, $false, $false, $true, 1, $false, "85-20=65", 2)
$d.Content.Find.Execute("34-13=21", $true, $false, $false, $false, $false, $true, 1, $false, "79+4=83", 2)
$d.Content.Find.Execute("93-43=50", $true, $false, $false, $false, $false, $true, 1, $false, "8-6=2", 2)
$d.Content.Find.Execute("96-71=25", $true, $false, $false, $false, $false, $true, 1, $false, "19+33=52", 2)
$d.Content.Find.Execute("5+9=14", $true, $false, $false, $false, $false, $true, 1, $false, "54-45=9", 2)
$d.Content.Find.Execute("83-11=72", $true, $false, $false, $false, $false, $true, 1, $false, "4+55=59", 2)
$d.Content.Find.Execute("2+1=3", $true, $false, $false, $false, $false, $true, 1, $false, "31+23=54", 2)
$d.Content.Find.Execute("65+28=93", $true, $false, $false, $false, $false, $true, 1, $false, "14+30=44", 2)
$d.Content.Find.Execute("75-53=22", $true, $false, $false, $false, $false, $true, 1, $false, "19+21=40", 2)
$d.Content.Find.Execute("32+65=97", $true, $false, $false, $false, $false, $true, 1, $false, "94-42=52", 2)
$d.Content.Find.Execute("2+33=35", $true, $false, $false, $false, $false, $true, 1, $false, "12+57=69", 2)
$d.Content.Find.Execute("10+83=93", $true, $false, $false, $false, $false, $true, 1, $false, "81-19=62", 2)
$d.Content.Find.Execute("92-18=74", $true, $false, $false, $false, $false, $true, 1, $false, "56-12=44", 2)
$d.Content.Find.Execute("29+13=42", $true, $false, $false, $false, $false, $true, 1, $false, "78-48=30", 2)
$d.Content.Find.Execute("60+14=74", $true, $false, $false, $false, $false, $true, 1, $false, "42-28=14", 2)
$d.Content.Find.Execute("8+56=64", $true, $false, $false, $false, $false, $true, 1, $false, "89-29=60", 2)
$d.Content.Find.Execute("96-83=13", $true, $false, $false, $false, $false, $true, 1, $false, "86+11=97", 2)
$d.Content.Find.Execute("47+42=89", $true, $false, $false, $false, $false, $true, 1, $false, "91+6=97", 2)
$d.Content.Find.Execute("10+1=11", $true, $false, $false, $false, $false, $true, 1, $false, "44-16=28", 2)
$d.Content.Find.Execute("74-1=73", $true, $false, $false, $false, $false, $true, 1, $false, "46-12=34", 2)
$d.Content.Find.Execute("96-54=42", $true, $false, $false, $false, $false, $true, 1, $false, "19+48=67", 2)
$d.Content.Find.Execute("8+15=23", $true, $false, $false, $false, $false, $true, 1, $false, "44-27=17", 2)
$d.Content.Find.Execute("97-27=70", $true, $false, $false, $false, $false, $true, 1, $false, "13-6=7", 2)
$d.Content.Find.Execute("39+59=98", $true, $false, $false, $false, $false, $true, 1, $false, "43-42=1", 2)
$d.Content.Find.Execute("77-38=39", $true, $false, $false, $false, $false, $true, 1, $false, "19+76=95", 2)
$d.Content.Find.Execute("84-78=6", $true, $false, $false, $false, $false, $true, 1, $false, "2+25=27", 2)
$d.Content.Find.Execute("38+34=72", $true, $false, $false, $false, $false, $true, 1, $false, "41+44=85", 2)
$d.Content.Find.Execute("7+34=41", $true, $false, $false, $false, $false, $true, 1, $false, "62-51=11", 2)
$d.Content.Find.Execute("32+26=58", $true, $false, $false, $false, $false, $true, 1, $false, "31-13=18", 2)
$d.Content.Find.Execute("35+53=88", $true, $false, $false, $false, $false, $true, 1, $false, "88-47=41", 2)
$d.Content.Find.Execute("34+24=58", $true, $false, $false, $false, $false, $true, 1, $false, "82-36=46", 2)
$d.Content.Find.Execute("53-21=32", $true, $false, $false, $false, $false, $true, 1, $false, "17+64=81", 2)
$d.Content.Find.Execute("49+28=77", $true, $false, $false, $false, $false, $true, 1, $false, "26+43=69", 2)
$d.Content.Find.Execute("96-11=85", $true, $false, $false, $false, $false, $true, 1, $false, "22-14=8", 2)
$d.Content.Find.Execute("45+35=80", $true, $false, $false, $false, $false, $true, 1, $false, "74-67=7", 2)
$d.Content.Find.Execute("43+28=71", $true, $false, $false, $false, $false, $true, 1, $false, "91-24=67", 2)
$d.Content.Find.Execute("6-4=2", $true, $false, $false, $false, $false, $true, 1, $false, "22-3=19", 2)
$d.Content.Find.Execute("48-5=43", $true, $false, $false, $false, $false, $true, 1, $false, "16+18=34", 2)
$d.Content.Find.Execute("36+8=44", $true, $false, $false, $false, $false, $true, 1, $false, "3+65=68", 2)
$d.Content.Find.Execute("1+33=34", $true, $false, $false, $false, $false, $true, 1, $false, "58+26=84", 2)
$d.Content.Find.Execute("63-52=11", $true, $false, $false, $false, $false, $true, 1, $false, "7+82=89", 2)
$d.Content.Find.Execute("68+8=76", $true, $false, $false, $false, $false, $true, 1, $false, "22-9=13", 2)
$d.Content.Find.Execute("81-59=22", $true, $false, $false, $false, $false, $true, 1, $false, "6+13=19", 2)
$d.Content.Find.Execute("7+27=34", $true, $false, $false, $false, $false, $true, 1, $false, "28+16=44", 2)
$d.Content.Find.Execute("64+14=78", $true, $false, $false, $false, $false, $true, 1, $false, "91-48=43", 2)
$d.Content.Find.Execute("78-74=4", $true, $false, $false, $false, $false, $true, 1, $false, "80-57=23", 2)
$d.Content.Find.Execute("83+2=85", $true, $false, $false, $false, $false, $true, 1, $false, "46+32=78", 2)
$d.Content.Find.Execute("53-36=17", $true, $false, $false, $false, $false, $true, 1, $false, "14+64=78", 2)
$d.Content.Find.Execute("19+60=79", $true, $false, $false, $false, $false, $true, 1, $false, "32+63=95", 2)
$d.Content.Find.Execute("32+2=34", $true, $false, $false, $false, $false, $true, 1, $false, "31-3=28", 2)
